# Weekly update of the "Agrícola del Norte S.A. de Arica - Ajo" price sheet.
# New weekly records are inserted; each new record bumps the older rows
# further down in the sheet (most-recent-first ordering), so we insert
# the new rows at their target positions and let Excel shift the rest
# of the table down automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the first new weekly record as row 35 ---------------------
$ws.Rows.Item(35).Insert()

$ws.Cells.Item(35,1).Value  = 1
$ws.Cells.Item(35,2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(35,3).Value  = "Arica y Parinacota"
$ws.Cells.Item(35,4).Value  = 44777
$ws.Cells.Item(35,5).Value  = 15
$ws.Cells.Item(35,6).Value  = 100112003
$ws.Cells.Item(35,7).Value  = "Ajo"
$ws.Cells.Item(35,8).Value  = "Chino"
$ws.Cells.Item(35,9).Value  = "Primera"
$ws.Cells.Item(35,10).Value = 200
$ws.Cells.Item(35,11).Value = 24000
$ws.Cells.Item(35,12).Value = 25000
$ws.Cells.Item(35,13).Value = 24500
$ws.Cells.Item(35,14).Value = "$/caja 10 kilos"
$ws.Cells.Item(35,15).Value = "China"
$ws.Cells.Item(35,16).Value = 2450
$ws.Cells.Item(35,17).Value = 10
$ws.Cells.Item(35,18).Value = "Hortaliza"

# --- Insert two more new weekly records as rows 39 and 40 -------------
$ws.Range("A39:A40").EntireRow.Insert()

$ws.Cells.Item(39,1).Value  = 1
$ws.Cells.Item(39,2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(39,3).Value  = "Arica y Parinacota"
$ws.Cells.Item(39,4).Value  = 45176
$ws.Cells.Item(39,5).Value  = 15
$ws.Cells.Item(39,6).Value  = 100112003
$ws.Cells.Item(39,7).Value  = "Ajo"
$ws.Cells.Item(39,8).Value  = "Chino"
$ws.Cells.Item(39,9).Value  = "Primera"
$ws.Cells.Item(39,10).Value = 300
$ws.Cells.Item(39,11).Value = 24000
$ws.Cells.Item(39,12).Value = 25000
$ws.Cells.Item(39,13).Value = 24500
$ws.Cells.Item(39,14).Value = "$/caja 10 kilos"
$ws.Cells.Item(39,15).Value = "China"
$ws.Cells.Item(39,16).Value = 2450
$ws.Cells.Item(39,17).Value = 10
$ws.Cells.Item(39,18).Value = "Hortaliza"

$ws.Cells.Item(40,1).Value  = 1
$ws.Cells.Item(40,2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(40,3).Value  = "Arica y Parinacota"
$ws.Cells.Item(40,4).Value  = 45176
$ws.Cells.Item(40,5).Value  = 15
$ws.Cells.Item(40,6).Value  = 100112003
$ws.Cells.Item(40,7).Value  = "Ajo"
$ws.Cells.Item(40,8).Value  = "Chino"
$ws.Cells.Item(40,9).Value  = "Segunda"
$ws.Cells.Item(40,10).Value = 400
$ws.Cells.Item(40,11).Value = 22000
$ws.Cells.Item(40,12).Value = 23000
$ws.Cells.Item(40,13).Value = 22500
$ws.Cells.Item(40,14).Value = "$/caja 10 kilos"
$ws.Cells.Item(40,15).Value = "China"
$ws.Cells.Item(40,16).Value = 2250
$ws.Cells.Item(40,17).Value = 10
$ws.Cells.Item(40,18).Value = "Hortaliza"
